$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Select a cell on Sheet1 to reflect new active selection, and zoom
$ws1.Range("B15").Select()
$excel.ActiveWindow.Zoom = 135

# Slightly adjusted column widths on Sheet1 (re-saved/re-measured widths)
$ws1.Columns.Item(1).ColumnWidth = 21.333333333333332
$ws1.Columns.Item(2).ColumnWidth = 26.333333333333332
$ws1.Columns.Item(3).ColumnWidth = 22.0

# Add the new "Questions" worksheet after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Questions"

# Populate header
$ws2.Range("A1").Value = "Qns for Elliot"

# Populate question column (entered out of strict row order, matching original authoring order)
$ws2.Range("A2").Value = "risk level?"
$ws2.Range("A4").Value = "source"
$ws2.Range("A5").Value = "lineitem compare @ price"
$ws2.Range("A3").Value = "if vendor =/= BBW?"
$ws2.Range("A6").Value = "lineitem discount"
$ws2.Range("A7").Value = "#6585 vs 6585 = same person?"
$ws2.Range("A8").Value = "wine virgin and wine nut priced @ 0?"

# Populate answer column
$ws2.Range("B2").Value = "security issues with payments. Ignore"
$ws2.Range("B3").Value = "this is the brand of the product"
$ws2.Range("B4").Value = "where the order comes from, no referrals?"
$ws2.Range("B5").Value = "comparison price = original price if the item is discounted"
$ws2.Range("B6").Value = "either % or real value"
$ws2.Range("B7").Value = "assume is the same"
$ws2.Range("B8").Value = "types of subscription. Usually items with 100% discount are given to subscribers - as the `$ is from the subscription renewa"

# Header cell formatting (bold) to match existing bold style used elsewhere in workbook
$ws2.Range("A1").Font.Bold = $true

# Column width + view state for the new sheet
$ws2.Columns.Item(1).ColumnWidth = 28.666666666666668

$ws2.Activate()
$ws2.Range("C6").Select()
$excel.ActiveWindow.Zoom = 133

Write-Host "done"
